# Updated Issue List for 1160, many issues added.  mem chart updated.
# Test files moved to Github and being added here.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Memory" sheet (sheet4): append rows 23-25 (memory map chart update)
# ---------------------------------------------------------------------------
$mem = $wb.Worksheets.Item("Memory")

# Row 23 - copy formatting from row 19 (A-only row, no B) for columns A,C:N
foreach ($col in @("A","C","D","E","F","G","H","I","J","K","L","M","N")) {
    $mem.Range($col + "19").Copy()
    $mem.Range($col + "23").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

$mem.Range("A23").Value = 1154
$mem.Range("C23").Value = "`$1FE6"
$mem.Range("D23").Value = "`$81E0"
$mem.Range("E23").Value = "`$BD00"
$mem.Range("F23").Value = "25082"
$mem.Range("G23").Value = "`$199F"
$mem.Range("H23").Value = "`$A200"
$mem.Range("I23").Value = "`$A200"
$mem.Range("J23").Value = "34913"
$mem.Range("K23").Value = "`$EC50"
$mem.Range("L23").Value = "`$EDBD"
$mem.Range("M23").Value = "`$FFFA"
$mem.Range("N23").Value = "04669"

# Row 24 - copy formatting from row 20 (B-only row, no A) for columns B,C:N
foreach ($col in @("B","C","D","E","F","G","H","I","J","K","L","M","N")) {
    $mem.Range($col + "20").Copy()
    $mem.Range($col + "24").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

$mem.Range("B24").Value = "1 tty"
$mem.Range("C24").Value = "`$1FE6"
$mem.Range("D24").Value = "`$87E0"
$mem.Range("E24").Value = "`$BD00"
$mem.Range("F24").Value = "26618"
$mem.Range("G24").Value = "`$199F"
$mem.Range("H24").Value = "`$A200"
$mem.Range("I24").Value = "`$A200"
$mem.Range("J24").Value = "34913"
$mem.Range("K24").Value = "`$EC50"
$mem.Range("L24").Value = "`$EDBD"
$mem.Range("M24").Value = "`$FFFA"
$mem.Range("N24").Value = "04669"

# Row 25 - copy formatting from row 19 again (A-only row, no B)
foreach ($col in @("A","C","D","E","F","G","H","I","J","K","L","M","N")) {
    $mem.Range($col + "19").Copy()
    $mem.Range($col + "25").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

$mem.Range("A25").Value = 1160
$mem.Range("C25").Value = "`$1FE6"
$mem.Range("D25").Value = "`$7CE0"
$mem.Range("E25").Value = "`$BD00"
$mem.Range("F25").Value = "23802"
$mem.Range("G25").Value = "`$199F"
$mem.Range("H25").Value = "`$A200"
$mem.Range("I25").Value = "`$A200"
$mem.Range("J25").Value = "34913"
$mem.Range("K25").Value = "`$EC50"
$mem.Range("L25").Value = "`$EDBD"
$mem.Range("M25").Value = "`$FFFA"
$mem.Range("N25").Value = "04669"

# ---------------------------------------------------------------------------
# 2) "Issues" sheet (sheet1): append issues 181-186 (rows 90-95)
# ---------------------------------------------------------------------------
$iss = $wb.Worksheets.Item("Issues")

function Add-Issue {
    param($row, $num, $category, $build, $desc, $height)

    foreach ($col in @("B","C","D","E")) {
        $iss.Range($col + "89").Copy()
        $iss.Range($col + $row).PasteSpecial(-4122)
    }
    $excel.CutCopyMode = $false

    $iss.Range("B$row").Value = $num
    $iss.Range("C$row").Value = $category
    $iss.Range("D$row").Value = $build
    $iss.Range("E$row").Value = $desc
    $iss.Rows.Item($row).RowHeight = $height
}

Add-Issue 90 181 "MORE" 1160 "When you more a file (more ../etc/init), if you are using the option of pressing return for one more line, when you get to the end of the file, return just gives you more and more blank lines instead of ending more because you are at the end of the file." 90

Add-Issue 91 182 "SHELL" 1160 "In a script, the CD command has no effect.  See TESTS/TESTS file which gets placed in root and should change you to usr/share/tests.  It does nothing.  Even if you do . TESTS.  I tried adding PWD after the CD command in the script, pwd shows that it the working dir changed, but script exits it goes back to wear it was." 105

Add-Issue 92 183 "SHELL" 1160 "At the end of TESTS/BUILDTEST (line 68) there is a dashed line ended with \n, which should on screen be a dashed line and then blank line, instead you get 2 dashed lines.  If you look at code, you can remove \n and the comment mark in line 69 and it displays right, but this second dashed line should not happen." 105

Add-Issue 93 184 "NETWORKING" 1160 "DHCP Never Leases and address.  If you do fixed IP address, DNS resolution does not appear to work and ping does not get any responses to local or remote networks.  TELNET locks after Ctrl-T and does not connect to any sources." 90

Add-Issue 94 185 "RETURN CODE" 1160 "RCTEST fails, you get bogus results.  Also in CDTEST, after each CD command you get some randowm number returned instead of 0 for no error (even though it worked).  Something is clearly off here." 75

Add-Issue 95 186 "MORE" 1160 "More Line numbering on long files works if you press space for next page, but repeats the same number over and over if you press return (though it does display contents right, except it will keep going passed end see 181)." 75

# ---------------------------------------------------------------------------
# 3) Active sheet / selection bookkeeping: Issues becomes the active tab,
#    Memory loses its "active" status, and each sheet's last selection moves
#    to reflect the newly-added rows.
# ---------------------------------------------------------------------------
$mem.Activate()
$mem.Range("N26").Select()

$iss.Activate()
$iss.Range("B96").Select()
